$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 27 ---
$ws.Cells.Item(27, 1).Value = 43769
$ws.Cells.Item(27, 2).Value = 0
$ws.Cells.Item(27, 3).Value = 0.083333333333333329
$ws.Cells.Item(27, 4).Value = 30
$ws.Cells.Item(27, 5).Value = 90
$t27 = "타팀 SRS Review"
$ws.Cells.Item(27, 6).Value = $t27
$c27 = $ws.Cells.Item(27, 6).Characters(3, $t27.Length - 2)
$c27.Font.Name = "돋움"
$c27.Font.Size = 10
$c27.Font.ColorIndex = -4105

# --- Row 28 ---
$ws.Cells.Item(28, 1).Value = 43770
$ws.Cells.Item(28, 2).Value = 0.83333333333333337
$ws.Cells.Item(28, 3).Value = 1
$ws.Cells.Item(28, 4).Value = 60
$ws.Cells.Item(28, 5).Value = 180
$t28 = "Skeleton Code refactoring - Function과 Controller 분리"
$ws.Cells.Item(28, 6).Value = $t28
$c28 = $ws.Cells.Item(28, 6).Characters(38, $t28.Length - 37)
$c28.Font.Name = "돋움"
$c28.Font.Size = 10
$c28.Font.ColorIndex = -4105

# --- Row 29 ---
$ws.Cells.Item(29, 1).Value = 43771
$ws.Cells.Item(29, 2).Value = 0.47916666666666669
$ws.Cells.Item(29, 3).Value = 0.54166666666666663
$ws.Cells.Item(29, 4).Value = 0
$ws.Cells.Item(29, 5).Value = 90
$ws.Cells.Item(29, 6).Value = $t28

# --- Row 30 ---
$ws.Cells.Item(30, 1).Value = 43770
$ws.Cells.Item(30, 2).Value = 0.66666666666666663
$ws.Cells.Item(30, 3).Value = 0.72916666666666663
$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(30, 5).Value = 90
$ws.Cells.Item(30, 6).Value = "Coding Guideline 작성"

# --- Row 31 ---
$ws.Cells.Item(31, 1).Value = 43775
$ws.Cells.Item(31, 2).Value = 0.75
$ws.Cells.Item(31, 3).Value = 0.83333333333333337
$ws.Cells.Item(31, 4).Value = 0
$ws.Cells.Item(31, 5).Value = 120
$t31 = "Coding Guideline 수정"
$ws.Cells.Item(31, 6).Value = $t31
$c31 = $ws.Cells.Item(31, 6).Characters(18, 2)
$c31.Font.Name = "돋움"
$c31.Font.Size = 10
$c31.Font.ColorIndex = -4105

# --- Sheet view: scroll position + selection ---
$ws.Range("B31").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
